# Update res_bus vm_pu results for the 380 kV case: slack bus voltage
# setpoint changed from 1.05 to 1.02 p.u., which re-solves the power flow
# and shifts every bus voltage magnitude (columns B-F and I-N) for rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$blockLeft = New-Object 'double[,]' 24,5
$blockLeft[0,0] = 1.02
$blockLeft[0,1] = 1.026490916865734
$blockLeft[0,2] = 1.025507051647289
$blockLeft[0,3] = 1.035281354641555
$blockLeft[0,4] = 1.043026249814855
$blockLeft[1,0] = 1.02
$blockLeft[1,1] = 1.027842156397847
$blockLeft[1,2] = 1.026019325364456
$blockLeft[1,3] = 1.03652365096261
$blockLeft[1,4] = 1.044448464980237
$blockLeft[2,0] = 1.02
$blockLeft[2,1] = 1.028716173043294
$blockLeft[2,2] = 1.026344013162027
$blockLeft[2,3] = 1.037327489195902
$blockLeft[2,4] = 1.045369012804245
$blockLeft[3,0] = 1.02
$blockLeft[3,1] = 1.02908353615165
$blockLeft[3,2] = 1.026478881440165
$blockLeft[3,3] = 1.037665423616462
$blockLeft[3,4] = 1.045756082122717
$blockLeft[4,0] = 1.02
$blockLeft[4,1] = 1.029145213854658
$blockLeft[4,2] = 1.026501430639359
$blockLeft[4,3] = 1.037722164438539
$blockLeft[4,4] = 1.045821077047701
$blockLeft[5,0] = 1.02
$blockLeft[5,1] = 1.028721082050322
$blockLeft[5,2] = 1.026345821692922
$blockLeft[5,3] = 1.037332004684755
$blockLeft[5,4] = 1.045374184559297
$blockLeft[6,0] = 1.02
$blockLeft[6,1] = 1.026947644726116
$blockLeft[6,2] = 1.025681579129091
$blockLeft[6,3] = 1.03570119870601
$blockLeft[6,4] = 1.043506838681758
$blockLeft[7,0] = 1.02
$blockLeft[7,1] = 1.023819926531252
$blockLeft[7,2] = 1.024459342914767
$blockLeft[7,3] = 1.032827278913259
$blockLeft[7,4] = 1.04021829376622
$blockLeft[8,0] = 1.02
$blockLeft[8,1] = 1.021732687406677
$blockLeft[8,2] = 1.023610008354521
$blockLeft[8,3] = 1.030910964937144
$blockLeft[8,4] = 1.038026988218159
$blockLeft[9,0] = 1.02
$blockLeft[9,1] = 1.020828330211808
$blockLeft[9,2] = 1.02323410637176
$blockLeft[9,3] = 1.030081045017767
$blockLeft[9,4] = 1.037078321730686
$blockLeft[10,0] = 1.02
$blockLeft[10,1] = 1.020492320724775
$blockLeft[10,2] = 1.023093261465912
$blockLeft[10,3] = 1.029772750100411
$blockLeft[10,4] = 1.036725967244162
$blockLeft[11,0] = 1.02
$blockLeft[11,1] = 1.020564400104801
$blockLeft[11,2] = 1.02312352824074
$blockLeft[11,3] = 1.029838881637243
$blockLeft[11,4] = 1.036801547541351
$blockLeft[12,0] = 1.02
$blockLeft[12,1] = 1.020800557450148
$blockLeft[12,2] = 1.023222488916466
$blockLeft[12,3] = 1.030055561812355
$blockLeft[12,4] = 1.037049195569475
$blockLeft[13,0] = 1.02
$blockLeft[13,1] = 1.020946049567315
$blockLeft[13,2] = 1.023283300548361
$blockLeft[13,3] = 1.030189062087951
$blockLeft[13,4] = 1.037201782521287
$blockLeft[14,0] = 1.02
$blockLeft[14,1] = 1.021792694161462
$blockLeft[14,2] = 1.023634784689578
$blockLeft[14,3] = 1.030966040613384
$blockLeft[14,4] = 1.038089951519511
$blockLeft[15,0] = 1.02
$blockLeft[15,1] = 1.022323615771931
$blockLeft[15,2] = 1.023853086440131
$blockLeft[15,3] = 1.031453377401597
$blockLeft[15,4] = 1.038647121976612
$blockLeft[16,0] = 1.02
$blockLeft[16,1] = 1.022633238848842
$blockLeft[16,2] = 1.02397963304905
$blockLeft[16,3] = 1.03173761916036
$blockLeft[16,4] = 1.038972128104004
$blockLeft[17,0] = 1.02
$blockLeft[17,1] = 1.02273880323527
$blockLeft[17,2] = 1.024022648879114
$blockLeft[17,3] = 1.031834536082859
$blockLeft[17,4] = 1.039082950004158
$blockLeft[18,0] = 1.02
$blockLeft[18,1] = 1.022266658587825
$blockLeft[18,2] = 1.023829745907683
$blockLeft[18,3] = 1.031401092208852
$blockLeft[18,4] = 1.038587341033059
$blockLeft[19,0] = 1.02
$blockLeft[19,1] = 1.020731017553108
$blockLeft[19,2] = 1.023193381061208
$blockLeft[19,3] = 1.029991755644657
$blockLeft[19,4] = 1.036976268817576
$blockLeft[20,0] = 1.02
$blockLeft[20,1] = 1.019764970093797
$blockLeft[20,2] = 1.022786227006253
$blockLeft[20,3] = 1.029105498996336
$blockLeft[20,4] = 1.035963451356113
$blockLeft[21,0] = 1.02
$blockLeft[21,1] = 1.020277142110416
$blockLeft[21,2] = 1.023002733729557
$blockLeft[21,3] = 1.029575336107435
$blockLeft[21,4] = 1.036500354652359
$blockLeft[22,0] = 1.02
$blockLeft[22,1] = 1.022292395255737
$blockLeft[22,2] = 1.023840294915894
$blockLeft[22,3] = 1.031424717678507
$blockLeft[22,4] = 1.038614353412237
$blockLeft[23,0] = 1.02
$blockLeft[23,1] = 1.024628864977666
$blockLeft[23,2] = 1.024781422877932
$blockLeft[23,3] = 1.033570305248743
$blockLeft[23,4] = 1.041068256048337
$ws.Range("B2:F25").Value = $blockLeft

$blockRight = New-Object 'double[,]' 24,6
$blockRight[0,0] = 1.028047310659273
$blockRight[0,1] = 1.031654064059758
$blockRight[0,2] = 1.028332420894459
$blockRight[0,3] = 1.03807840888099
$blockRight[0,4] = 1.045801277101354
$blockRight[0,5] = 1.033119131163119
$blockRight[1,0] = 1.028114196249099
$blockRight[1,1] = 1.032643018763664
$blockRight[1,2] = 1.02865281339994
$blockRight[1,3] = 1.039128916684577
$blockRight[1,4] = 1.047032844707646
$blockRight[1,5] = 1.034109490296118
$blockRight[2,0] = 1.028152502439901
$blockRight[2,1] = 1.033282201577143
$blockRight[2,2] = 1.028852644648525
$blockRight[2,3] = 1.039808101438636
$blockRight[2,4] = 1.047829512859632
$blockRight[2,5] = 1.034749580822508
$blockRight[3,0] = 1.028167412619882
$blockRight[3,1] = 1.033550740576031
$blockRight[3,2] = 1.028934856412712
$blockRight[3,3] = 1.040093498731202
$blockRight[3,4] = 1.048164378682872
$blockRight[3,5] = 1.035018501177576
$blockRight[4,0] = 1.0281698460087
$blockRight[4,1] = 1.033595819370691
$blockRight[4,2] = 1.028948554567337
$blockRight[4,3] = 1.040141410541831
$blockRight[4,4] = 1.048220601097659
$blockRight[4,5] = 1.035063643989295
$blockRight[5,0] = 1.028152706365405
$blockRight[5,1] = 1.033285790486255
$blockRight[5,2] = 1.028853750236068
$blockRight[5,3] = 1.039811915446027
$blockRight[5,4] = 1.047833987558758
$blockRight[5,5] = 1.034753174828282
$blockRight[6,0] = 1.02807094286734
$blockRight[6,1] = 1.031988440333347
$blockRight[6,2] = 1.028442245416833
$blockRight[6,3] = 1.038433551539637
$blockRight[6,4] = 1.046217542688354
$blockRight[6,5] = 1.033453982289349
$blockRight[7,0] = 1.027888906352644
$blockRight[7,1] = 1.029696565096052
$blockRight[7,2] = 1.027660033183786
$blockRight[7,3] = 1.036000241370182
$blockRight[7,4] = 1.043367172052731
$blockRight[7,5] = 1.031158852326335
$blockRight[8,0] = 1.027742192411911
$blockRight[8,1] = 1.028164574574658
$blockRight[8,2] = 1.027100466945094
$blockRight[8,3] = 1.034374845117876
$blockRight[8,4] = 1.041465360200861
$blockRight[8,5] = 1.029624686202679
$blockRight[9,0] = 1.027672681232815
$blockRight[9,1] = 1.027500197155752
$blockRight[9,2] = 1.026849189472056
$blockRight[9,3] = 1.033670230717229
$blockRight[9,4] = 1.040641432365732
$blockRight[9,5] = 1.028959365291633
$blockRight[10,0] = 1.027645965106361
$blockRight[10,1] = 1.027253261528089
$blockRight[10,2] = 1.026754508480939
$blockRight[10,3] = 1.033408380493541
$blockRight[10,4] = 1.040335319692193
$blockRight[10,5] = 1.028712078987055
$blockRight[11,0] = 1.027651736342748
$blockRight[11,1] = 1.027306237190017
$blockRight[11,2] = 1.02677487869296
$blockRight[11,3] = 1.033464553975997
$blockRight[11,4] = 1.040400985102832
$blockRight[11,5] = 1.028765129880499
$blockRight[12,0] = 1.027670491146611
$blockRight[12,1] = 1.027479788589222
$blockRight[12,2] = 1.026841390534216
$blockRight[12,3] = 1.033648588669355
$blockRight[12,4] = 1.04061613043405
$blockRight[12,5] = 1.028938927742598
$blockRight[13,0] = 1.027681927861248
$blockRight[13,1] = 1.027586698515009
$blockRight[13,2] = 1.026882192518695
$blockRight[13,3] = 1.033761961831765
$blockRight[13,4] = 1.040748679270769
$blockRight[13,5] = 1.029045989492742
$blockRight[14,0] = 1.027746679784026
$blockRight[14,1] = 1.028208645405253
$blockRight[14,2] = 1.02711695452225
$blockRight[14,3] = 1.034421590726031
$blockRight[14,4] = 1.041520032120772
$blockRight[14,5] = 1.029668819618909
$blockRight[15,0] = 1.027785696793537
$blockRight[15,1] = 1.028598501782261
$blockRight[15,2] = 1.027261812832631
$blockRight[15,3] = 1.034835139739719
$blockRight[15,4] = 1.042003762906309
$blockRight[15,5] = 1.030059229636689
$blockRight[16,0] = 1.027807877329898
$blockRight[16,1] = 1.02882580069113
$blockRight[16,2] = 1.027345439387844
$blockRight[16,3] = 1.035076278268112
$blockRight[16,4] = 1.042285873631966
$blockRight[16,5] = 1.03028685133608
$blockRight[17,0] = 1.027815342305294
$blockRight[17,1] = 1.02890328728206
$blockRight[17,2] = 1.027373806722037
$blockRight[17,3] = 1.035158487224475
$blockRight[17,4] = 1.042382059235768
$blockRight[17,5] = 1.030364447966857
$blockRight[18,0] = 1.027781570346936
$blockRight[18,1] = 1.02855668401827
$blockRight[18,2] = 1.027246360540018
$blockRight[18,3] = 1.034790777901108
$blockRight[18,4] = 1.041951867472192
$blockRight[18,5] = 1.030017352486677
$blockRight[19,0] = 1.027664993056495
$blockRight[19,1] = 1.02742868637165
$blockRight[19,2] = 1.026821841561943
$blockRight[19,3] = 1.033594398522078
$blockRight[19,4] = 1.040552777442377
$blockRight[19,5] = 1.028887752954016
$blockRight[20,0] = 1.027586509990163
$blockRight[20,1] = 1.026718563481539
$blockRight[20,2] = 1.026547147649343
$blockRight[20,3] = 1.032841461854167
$blockRight[20,4] = 1.039672711468772
$blockRight[20,5] = 1.028176621607948
$blockRight[21,0] = 1.027628606217047
$blockRight[21,1] = 1.027095100184892
$blockRight[21,2] = 1.026693504457306
$blockRight[21,3] = 1.033240677867888
$blockRight[21,4] = 1.040139290747095
$blockRight[21,5] = 1.028553693036609
$blockRight[22,0] = 1.027783436694908
$blockRight[22,1] = 1.028575579967526
$blockRight[22,2] = 1.027253345444725
$blockRight[22,3] = 1.034810823346057
$blockRight[22,4] = 1.04197531691229
$blockRight[22,5] = 1.030036275270348
$blockRight[23,0] = 1.027940447854658
$blockRight[23,1] = 1.030289773934296
$blockRight[23,2] = 1.027868987877461
$blockRight[23,3] = 1.036629858356135
$blockRight[23,4] = 1.044104320226557
$blockRight[23,5] = 1.031752903589171
$ws.Range("I2:N25").Value = $blockRight
